$d = $word.ActiveDocument

$replacements = @(
    @("953÷5=", "235÷9="),
    @("849÷9=", "780÷6="),
    @("641÷9=", "469÷9="),
    @("458÷7=", "828÷3="),
    @("582÷4=", "924÷8="),
    @("110÷6=", "539÷8="),
    @("214÷8=", "823÷6="),
    @("611÷3=", "896÷9="),
    @("751÷4=", "105÷5="),
    @("180÷6=", "214÷6="),
    @("618÷6=", "312÷8="),
    @("152÷5=", "781÷2="),
    @("340÷2=", "879÷8="),
    @("356÷6=", "123÷3="),
    @("185÷4=", "291÷6="),
    @("142÷3=", "120÷6="),
    @("719÷9=", "153÷3="),
    @("901÷2=", "463÷3="),
    @("123÷7=", "670÷5="),
    @("988÷3=", "341÷7="),
    @("416÷6=", "899÷5="),
    @("904÷7=", "887÷2="),
    @("163÷2=", "796÷5="),
    @("512÷7=", "846÷4="),
    @("219÷9=", "299÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
